# Apply the "fabrics table" refactor to the DB Management workbook.
#
# The sheet holds three mini "tables" drawn side by side:
#   B:D  -> "Style Group" table
#   F:H  -> "POs" table
#   J:L  -> "colors" table
#
# The commit adds a fabrics table to the DB design and renames a few
# fields to make the schema consistent. Concretely (relative to the
# current values already on the sheet):
#   J2 "team/color"   -> "colors"       (table title)
#   K3 "team"         -> "color_code"
#   K4 "color_combo"  -> "team"
#   C6 "Paking Method"-> "Type"
# and the active selection moves from K6 to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "colors"
$ws.Range("K3").Value = "color_code"
$ws.Range("K4").Value = "team"
$ws.Range("C6").Value = "Type"

$ws.Range("C7").Select()
